$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the rnaDate column (D): was mistakenly "S.GISH" (a name), should be the date "05.25.10".
# Force text entry (no number-format styling) so it stays a plain shared string like before.
$rngD = $ws.Range("D2:D5")
$rngD.NumberFormat = "@"
$rngD.Value = "05.25.10"
$rngD.ClearFormats()

# Correct the rnaPreparer column (E): was mistakenly "Retrofitted_159", should be "S.GISH".
$ws.Range("E2:E5").Value = "S.GISH"
